
# Generate Report for Archive
# Update the localization status from "Ready for handoff" to "In Translation"
# across the Overview sheet (zh-cn / de-de status columns) and the two
# per-locale detail sheets (their "Status" column), then resize the
# now-narrower "Status" columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Shrink the Status-related columns to match the shorter text.
# (12.5 "characters" is the input that Excel's pixel-grid rounding maps
# closest to the target stored width of ~13.41.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
